$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D, L, M, N, O, P, Q, R, S, T
# (the underlying weekly data rows were reordered/updated)
$data = @{
    2 = @{ D = 44742; L = "Segunda"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 806;   T = 18 }
    3 = @{ D = 44330; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó";  S = 861;   T = 18 }
    4 = @{ D = 44719; L = "Primera"; M = 50;  N = 14000; O = 15000; P = 14400; Q = "`$/caja 18 kilos granel"; R = "Región del Maule";       S = 800;   T = 18 }
    5 = @{ D = 44707; L = "Primera"; M = 60;  N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Provincia de Curicó";  S = 1042;  T = 12 }
    6 = @{ D = 44334; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 12 kilos granel"; R = "Región de O'Higgins";    S = 11500; T = 1 }
    7 = @{ D = 44708; L = "Primera"; M = 70;  N = 12000; O = 13000; P = 12571; Q = "`$/caja 12 kilos empedrada"; R = "Provincia de Curicó";  S = 1048;  T = 12 }
    8 = @{ D = 44714; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins";     S = 806;   T = 18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
